# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking Price values are written with a leading apostrophe (forces
# text entry, like typing into the cell in Excel) and then restyled back to
# "Normal" so they stay plain text (matching the sheet's inlineStr cells)
# without picking up a stray quote-prefix cell style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.849.88"
$ws.Range("D3").Value = "3.096.97"
$ws.Range("E3").Value = "  +5.20%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'580.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").Value = "'172.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.62%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.090.53"
$ws.Range("E8").Value = "  +5.11%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'6.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.73%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.85%  "
$ws.Range("E12").Value = "  +4.03%  "
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").Value = "'37.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.61%  "
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "3.610.15"
$ws.Range("E16").Value = "  +5.11%  "
$ws.Range("D17").Value = "66.843.72"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").Value = "'7.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "3.100.56"
$ws.Range("E19").Value = "  +5.27%  "
$ws.Range("D20").Value = "'16.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").Value = "'479.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.50%  "
$ws.Range("D22").Value = "'0.714"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("E23").Value = "  +3.16%  "
$ws.Range("D24").Value = "'83.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.18%  "
$ws.Range("E25").Value = "  +8.44%  "
$ws.Range("D26").Value = "'2.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.63%  "
$ws.Range("D27").Value = "'10.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'7.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("E30").Value = "  -2.59%  "
$ws.Range("D31").Value = "'2.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("D32").Value = "'28.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.74%  "
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").Value = "'0.114"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'5.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("D37").Value = "'0.990"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("D38").Value = "'48.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("E39").Value = "  +6.97%  "
$ws.Range("D40").Value = "'50.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("D41").Value = "'0.315"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.24%  "
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").Value = "'8.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").Value = "'2.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("D45").Value = "2.839.65"
$ws.Range("E45").Value = "  +6.24%  "
$ws.Range("D46").Value = "'0.0360"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("D47").Value = "'384.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'135.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'24.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("E51").Value = "  +3.00%  "
